# syneos_word_template.docx edit
#
# 1. The (empty) first body paragraph switches from direct character
#    formatting (sz=22) to the built-in "Heading1" paragraph style.
# 2. The "Heading1" / "Heading1Char" styles switch their font colour
#    from the theme's Background 1 (white) to Text 1 (near-black) and
#    pick up an explicit en-GB language tag.

$d = $word.ActiveDocument

# --- 1. First paragraph becomes a Heading 1 -------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Style = "Heading1"

# --- 2. Re-colour the Heading 1 style (and its linked character style) ----
# wdThemeColorText1 = 13 (was wdThemeColorBackground1 = 12)
$wdThemeColorText1 = 13

$headingStyle = $d.Styles("Heading1")
$headingStyle.Font.TextColor.ObjectThemeColor = $wdThemeColorText1
$headingStyle.LanguageID = "en-GB"

$headingCharStyle = $d.Styles("Heading1Char")
$headingCharStyle.Font.TextColor.ObjectThemeColor = $wdThemeColorText1
$headingCharStyle.LanguageID = "en-GB"
